# The "Förändrad" (Changed) column C for data rows 2-374 is updated from
# date serial 45205 (2023-10-06) to 45206 (2023-10-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = 374

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 45206
    }
}
